$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns (row 2) for the DATASET_TYPE table
$ws.Range("D2").Value = "Main Data Set Pattern"
$ws.Range("E2").Value = "Main Data Set Path"
$ws.Range("F2").Value = "Disallow Deletion"
$ws.Range("G2").Value = "Modification Date"

# New data values (row 3) for the ATTACHMENT data set type
$ws.Range("D3").Value = ".*\.jpg"
$ws.Range("E3").Value = "original/images/"

# "FALSE" must land as literal text (matching the existing FALSE text cells
# elsewhere in the sheet), not as a native boolean - copy an existing text
# "FALSE" cell's value across via PasteSpecial(values only) so no boolean
# coercion happens and no extra formatting is carried over.
$ws.Range("B5").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4163) | Out-Null

$ws.Range("G3").Value = "2023-03-10 17:23:44"
$ws.Range("G3").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

$ws.Range("G4").Select()
